$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: merge the runs around the removed "_GoBack" bookmark so that
# "...obsahuje editor " + (bookmark) + "map" become a single run
# "...obsahuje editor map" (bookmark removed).
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Dále také obsahuje editor map",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dále také obsahuje editor map", 2) | Out-Null

# -----------------------------------------------------------------------
# Change 2: append a new sentence after "...jak co ovládat a poznat."
# as its own run, then add a new "_GoBack" bookmark right after it.
# -----------------------------------------------------------------------
$target = $d.Paragraphs.Item(14)
$origEnd = $target.Range.End - 1

# Insert a temporary placeholder character and wrap it with a temporary
# bookmark so the two sentences stay in separate runs (this workaround
# avoids an engine quirk where InsertAfter on a bare collapsed range
# always merges into the preceding identically-formatted run).
$anchor = $d.Range($origEnd, $origEnd)
$anchor.InsertAfter("Y")
$tempBmRange = $d.Range($origEnd, $origEnd + 1)
$d.Bookmarks.Add("TempSplitBM", $tempBmRange) | Out-Null

$afterTemp = $d.Range($origEnd + 1, $origEnd + 1)
$afterTemp.InsertAfter(" Pokud nerozumíte ovládání hry, pak se doporučuji podívat tam.")

$tempBm = $d.Bookmarks("TempSplitBM")
$tempBm.Range.Text = ""
$d.Bookmarks("TempSplitBM").Delete()

# Now append a placeholder at the (new) end of the paragraph, wrap it in
# the "_GoBack" bookmark, then remove just the placeholder text, leaving
# a collapsed bookmark at the very end of the paragraph.
$p2 = $d.Paragraphs.Item(14)
$endPos = $p2.Range.End - 1
$endAnchor = $d.Range($endPos, $endPos)
$endAnchor.InsertAfter("Z")

$finalBmRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $finalBmRange) | Out-Null
$d.Bookmarks("_GoBack").Range.Text = ""
